# Fix inconsistent capitalization in "Years collected" and "Ageing structure"
# columns on the single worksheet, matching already-used canonical spellings
# elsewhere in the workbook (e.g. "1983 To 2022" -> "1983 to 2022",
# "1970 - Present" -> "1970 - present", "Shell" -> "shell", "Scale" -> "scale").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I15").Value = "1983 to 2022"
$ws.Range("I16").Value = "1973 to 2022"

$ws.Range("I19").Value = "1970 - present"
$ws.Range("I20").Value = "1970 - present"
$ws.Range("I21").Value = "1970 - present"
$ws.Range("I22").Value = "1970 - present"

$ws.Range("F25").Value = "shell"
$ws.Range("F26").Value = "shell"

$ws.Range("F27").Value = "scale"
$ws.Range("F28").Value = "scale"
$ws.Range("F29").Value = "scale"

$wb.Save()
